$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (existing row, values replaced) ---
$ws.Range("A2").Value = "aaa"
$ws.Range("B2").Value = "aaa"
$ws.Range("C2").Value = "aaa"
$ws.Range("D2").Value = 111
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 11.1
$ws.Range("G2").Value = 1500
$ws.Range("H2").Value = "18/07/2024"
$ws.Range("I2").Value = "18/07/2025"
$ws.Range("J2").Value = $false
# K2 / L2 stay as they already were (empty) - untouched

# --- Row 3 (new row) ---
$ws.Range("A3").Value = "bbb"
$ws.Range("B3").Value = "aaa"
$ws.Range("C3").Value = "aaa"
$ws.Range("D3").Value = 15000
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 2250
$ws.Range("G3").Value = 5000
$ws.Range("H3").Value = "18/02/2024"
$ws.Range("I3").Value = "18/02/2025"
$ws.Range("J3").Value = $false

# --- Row 4 (new row) ---
$ws.Range("A4").Value = "ccc"
$ws.Range("B4").Value = "bbb"
$ws.Range("C4").Value = "ccc"
$ws.Range("D4").Value = 1500
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 150
$ws.Range("G4").Value = 155
$ws.Range("H4").Value = "15/07/2004"
$ws.Range("I4").Value = "15/07/2005"
$ws.Range("J4").Value = $false

# --- Row 5 (new row) ---
$ws.Range("A5").Value = "Hugo Rios Brito"
$ws.Range("B5").Value = "Porsche GT3"
$ws.Range("C5").Value = "BRADESCO SEGUROS"

# D5..G5 and L5 are stored as text (not numbers) in the target even though
# they look like numbers. A leading apostrophe forces Excel to keep the
# literal text (same trick used interactively); ClearFormats() afterwards
# drops the leftover "quote prefix" style so the cell keeps the default
# (unstyled) look, matching the target.
$ws.Range("D5").Value = "'20000"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'10"
$ws.Range("E5").ClearFormats()
$ws.Range("F5").Value = "'2000.00"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").Value = "'3500"
$ws.Range("G5").ClearFormats()

$ws.Range("H5").Value = "15/07/2004"
$ws.Range("I5").Value = "15/07/2005"
$ws.Range("J5").Value = $true
$ws.Range("K5").Value = "THEUREN"

$ws.Range("L5").Value = "'500"
$ws.Range("L5").ClearFormats()
